$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a serial date value that was bumped by one day
# (45179 -> 45180) for every data row (rows 2 through 151).
$ws.Range("C2:C151").Value = 45180
